# Fixed proc Use and Modify
# The "Line number" column (A) for the generated test rows had an off-by-one
# gap (row 152 jumped from 93 to 95, skipping 94) for every row from 152
# through 226. Decrement each of those line numbers by 1 so the sequence is
# continuous again (1..142).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 152; $r -le 226; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value2 = $val - 1
    }
}

# Restore the view to where the author last left it while reviewing the fix.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 207
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("R223").Select() | Out-Null
